# Serbia Super Liga - odds-base update (04-04-2024 23:22)
#
# The source feed re-synced a handful of fixture rows:
#   - rows 195/196, 212/213 and 214/215 each had their data swapped
#     (columns B..AC; the running "id" in column A stays put)
#   - row 221 (fixture 6979598, not yet played) was dropped entirely
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    # Swap every data column (B..AC) between two rows, leaving column A
    # (the positional id) untouched. Use Value2 (not Value) so the
    # full double precision of the odds survives the round trip.
    $rangeA = $ws.Range("B${rowA}:AC${rowA}")
    $rangeB = $ws.Range("B${rowB}:AC${rowB}")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-Rows 195 196
Swap-Rows 212 213
Swap-Rows 214 215

# Drop the still-unplayed fixture that was pulled from the feed; this
# shifts nothing else since it was the last data row (220 -> new last row).
$ws.Rows(221).Delete()
